$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.131.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4710"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3684"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8810"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.849.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07313"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.462"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.549"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008766"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.141.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.310"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.069.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.898"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.158"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.272"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7594"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.169"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.534"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.930"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05332"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.431"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1660"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.546"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4945"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06315"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
